$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 (C naive3b) - update to consistent baseline for approach 3b
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 36
$ws.Range("O6").Value = 7.319620253164557
$ws.Range("P6").Value = 8.506340931670369
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 44
$ws.Range("V6").Value = 2.810126582278481
$ws.Range("W6").Value = 4.972498204788462
$ws.Range("AB6").Value = 33
$ws.Range("AC6").Value = 1.338607594936709
$ws.Range("AD6").Value = 4.130709223604813
$ws.Range("AJ6").Value = 0.129746835443038
$ws.Range("AK6").Value = 0.3637563275916581
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 1.473684210526316
$ws.Range("AO6").Value = 2.5
$ws.Range("AP6").Value = 5
$ws.Range("AQ6").Value = 1.949513816186035
$ws.Range("AR6").Value = 0.712827936278364
$ws.Range("AS6").Value = 11.56756756756757
$ws.Range("AT6").Value = 15.8125
$ws.Range("AU6").Value = 17.41666666666666
$ws.Range("AV6").Value = 19.4
$ws.Range("AW6").Value = 59
$ws.Range("AX6").Value = 18.53251213318498
$ws.Range("AY6").Value = 5.580829287170561

# Row 11 (P naive3b) - update to consistent baseline for approach 3b
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 6
$ws.Range("N11").Value = 34
$ws.Range("O11").Value = 6.55511811023622
$ws.Range("P11").Value = 5.867687475497098
$ws.Range("T11").Value = 4
$ws.Range("U11").Value = 21
$ws.Range("V11").Value = 3.204724409448819
$ws.Range("W11").Value = 3.530417979435353
$ws.Range("AB11").Value = 26
$ws.Range("AC11").Value = 2.122047244094488
$ws.Range("AD11").Value = 4.446083027817777
$ws.Range("AJ11").Value = 0.2322834645669291
$ws.Range("AK11").Value = 0.4758818549636049
$ws.Range("AM11").Value = 1.157894736842105
$ws.Range("AO11").Value = 2.5
$ws.Range("AP11").Value = 4
$ws.Range("AQ11").Value = 1.812612239259566
$ws.Range("AR11").Value = 0.858066687701709
$ws.Range("AS11").Value = 10.55555555555556
$ws.Range("AT11").Value = 15.83333333333333
$ws.Range("AU11").Value = 19
$ws.Range("AV11").Value = 23.75
$ws.Range("AW11").Value = 95
$ws.Range("AX11").Value = 22.00835514543718
$ws.Range("AY11").Value = 13.28140617026885
